# Auto-generated edit script: update cached market-price / profit figures
# across the Leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect a refreshed data pull ("chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 4 (Leve Item ID 5470)
$ws.Cells.Item(4, 8).Value = 160
$ws.Cells.Item(4, 9).Value = 160
$ws.Cells.Item(4, 11).Value = 160
$ws.Cells.Item(4, 13).Value = -46
# row 9 (Leve Item ID 5487)
$ws.Cells.Item(9, 8).Value = 182.66667
$ws.Cells.Item(9, 9).Value = 111.5
$ws.Cells.Item(9, 10).Value = 325
$ws.Cells.Item(9, 11).Value = 111.5
$ws.Cells.Item(9, 12).Value = 325
$ws.Cells.Item(9, 13).Value = 57.5
$ws.Cells.Item(9, 14).Value = -663
# row 18 (Leve Item ID 5471)
$ws.Cells.Item(18, 8).Value = 4711.9375
$ws.Cells.Item(18, 9).Value = 2900
$ws.Cells.Item(18, 10).Value = 6523.875
$ws.Cells.Item(18, 11).Value = 2900
$ws.Cells.Item(18, 12).Value = 6523.875
$ws.Cells.Item(18, 13).Value = -2616
$ws.Cells.Item(18, 14).Value = -7091.875
# row 33 (Leve Item ID 5512)
$ws.Cells.Item(33, 8).Value = 226.25
$ws.Cells.Item(33, 9).Value = 226.25
$ws.Cells.Item(33, 11).Value = 226.25
$ws.Cells.Item(33, 13).Value = 2.75
# row 34 (Leve Item ID 2160)
$ws.Cells.Item(34, 8).Value = 2932.3333
$ws.Cells.Item(34, 9).Value = 2508.8
$ws.Cells.Item(34, 10).Value = 5050
$ws.Cells.Item(34, 11).Value = 2508.8
$ws.Cells.Item(34, 12).Value = 5050
$ws.Cells.Item(34, 13).Value = -2305.8
$ws.Cells.Item(34, 14).Value = -5456
# row 36 (Leve Item ID 2160)
$ws.Cells.Item(36, 8).Value = 2932.3333
$ws.Cells.Item(36, 9).Value = 2508.8
$ws.Cells.Item(36, 10).Value = 5050
$ws.Cells.Item(36, 11).Value = 2508.8
$ws.Cells.Item(36, 12).Value = 5050
$ws.Cells.Item(36, 13).Value = -1793.8
$ws.Cells.Item(36, 14).Value = -6480
# row 58 (Leve Item ID 4606)
$ws.Cells.Item(58, 8).Value = 4017.75
$ws.Cells.Item(58, 10).Value = 7857.143
$ws.Cells.Item(58, 12).Value = 23571.429
$ws.Cells.Item(58, 14).Value = -23871.429
# row 93 (Leve Item ID 18043)
$ws.Cells.Item(93, 8).Value = 64999
$ws.Cells.Item(93, 10).Value = 64999
$ws.Cells.Item(93, 12).Value = 64999
$ws.Cells.Item(93, 14).Value = -69991
# row 98 (Leve Item ID 36237)
$ws.Cells.Item(98, 8).Value = 2338.92
$ws.Cells.Item(98, 9).Value = 1227.5333
$ws.Cells.Item(98, 11).Value = 1227.5333
$ws.Cells.Item(98, 13).Value = 270.4666999999999
# row 106 (Leve Item ID 19903)
$ws.Cells.Item(106, 8).Value = 62521372
$ws.Cells.Item(106, 10).Value = 24570.572
$ws.Cells.Item(106, 12).Value = 24570.572
$ws.Cells.Item(106, 14).Value = -25832.572
# row 107 (Leve Item ID 27766)
$ws.Cells.Item(107, 8).Value = 45455548
$ws.Cells.Item(107, 9).Value = 50001052
$ws.Cells.Item(107, 11).Value = 50001052
$ws.Cells.Item(107, 13).Value = -49999132
# row 122 (Leve Item ID 36237)
$ws.Cells.Item(122, 8).Value = 2338.92
$ws.Cells.Item(122, 9).Value = 1227.5333
$ws.Cells.Item(122, 11).Value = 3682.5999
$ws.Cells.Item(122, 13).Value = -1232.5999

$ws = $wb.Worksheets.Item("ARM")
# row 55 (Leve Item ID 2830)
$ws.Cells.Item(55, 8).Value = 19474.334
$ws.Cells.Item(55, 9).Value = 4848
$ws.Cells.Item(55, 10).Value = 22399.6
$ws.Cells.Item(55, 11).Value = 4848
$ws.Cells.Item(55, 12).Value = 22399.6
$ws.Cells.Item(55, 13).Value = -4533
$ws.Cells.Item(55, 14).Value = -23029.6

$ws = $wb.Worksheets.Item("BSM")
# row 86 (Leve Item ID 12526)
$ws.Cells.Item(86, 8).Value = 3441.6
$ws.Cells.Item(86, 9).Value = 3329.5454
$ws.Cells.Item(86, 10).Value = 3749.75
$ws.Cells.Item(86, 11).Value = 3329.5454
$ws.Cells.Item(86, 12).Value = 3749.75
$ws.Cells.Item(86, 13).Value = -2206.5454
$ws.Cells.Item(86, 14).Value = -5995.75
# row 89 (Leve Item ID 12526)
$ws.Cells.Item(89, 8).Value = 3441.6
$ws.Cells.Item(89, 9).Value = 3329.5454
$ws.Cells.Item(89, 10).Value = 3749.75
$ws.Cells.Item(89, 11).Value = 16647.727
$ws.Cells.Item(89, 12).Value = 18748.75
$ws.Cells.Item(89, 13).Value = -11031.727
$ws.Cells.Item(89, 14).Value = -29980.75
# row 99 (Leve Item ID 19943)
$ws.Cells.Item(99, 8).Value = 9749.5
$ws.Cells.Item(99, 9).Value = 3000
$ws.Cells.Item(99, 11).Value = 3000
$ws.Cells.Item(99, 13).Value = -1502
# row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 3110.6667
$ws.Cells.Item(134, 9).Value = 2999
$ws.Cells.Item(134, 11).Value = 8997
$ws.Cells.Item(134, 13).Value = -6462

$ws = $wb.Worksheets.Item("CRP")
# row 22 (Leve Item ID 5367)
$ws.Cells.Item(22, 8).Value = 86177.21000000001
$ws.Cells.Item(22, 9).Value = 100373.836
$ws.Cells.Item(22, 11).Value = 100373.836
$ws.Cells.Item(22, 13).Value = -100023.836
# row 74 (Leve Item ID 10636)
$ws.Cells.Item(74, 8).Value = 59999
$ws.Cells.Item(74, 10).Value = 59999
$ws.Cells.Item(74, 12).Value = 59999
$ws.Cells.Item(74, 14).Value = -61747
# row 77 (Leve Item ID 10636)
$ws.Cells.Item(77, 8).Value = 59999
$ws.Cells.Item(77, 10).Value = 59999
$ws.Cells.Item(77, 12).Value = 179997
$ws.Cells.Item(77, 14).Value = -188733
# row 103 (Leve Item ID 19558)
$ws.Cells.Item(103, 8).Value = 37302.855
$ws.Cells.Item(103, 9).Value = 27424.2
$ws.Cells.Item(103, 11).Value = 27424.2
$ws.Cells.Item(103, 13).Value = -26252.2

$ws = $wb.Worksheets.Item("CUL")
# row 2 (Leve Item ID 4847)
$ws.Cells.Item(2, 8).Value = 39.434784
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 12).Value = 480
$ws.Cells.Item(2, 14).Value = -706
# row 4 (Leve Item ID 4650)
$ws.Cells.Item(4, 8).Value = 2007697.8
$ws.Cells.Item(4, 9).Value = 2215982.8
$ws.Cells.Item(4, 11).Value = 6647948.399999999
$ws.Cells.Item(4, 13).Value = -6647836.399999999
# row 137 (Leve Item ID 44088)
$ws.Cells.Item(137, 8).Value = 4210.3335
$ws.Cells.Item(137, 9).Value = 4330.5
$ws.Cells.Item(137, 10).Value = 3970
$ws.Cells.Item(137, 11).Value = 12991.5
$ws.Cells.Item(137, 12).Value = 11910
$ws.Cells.Item(137, 13).Value = -7891.5
$ws.Cells.Item(137, 14).Value = -22110

$ws = $wb.Worksheets.Item("GSM")
# row 80 (Leve Item ID 12521)
$ws.Cells.Item(80, 8).Value = 9130.799999999999
$ws.Cells.Item(80, 10).Value = 11351
$ws.Cells.Item(80, 12).Value = 11351
$ws.Cells.Item(80, 14).Value = -13347
# row 83 (Leve Item ID 12521)
$ws.Cells.Item(83, 8).Value = 9130.799999999999
$ws.Cells.Item(83, 10).Value = 11351
$ws.Cells.Item(83, 12).Value = 56755
$ws.Cells.Item(83, 14).Value = -66739
# row 107 (Leve Item ID 27802)
$ws.Cells.Item(107, 8).Value = 1125.5454
$ws.Cells.Item(107, 9).Value = 869.1429000000001
$ws.Cells.Item(107, 10).Value = 1574.25
$ws.Cells.Item(107, 11).Value = 869.1429000000001
$ws.Cells.Item(107, 12).Value = 1574.25
$ws.Cells.Item(107, 13).Value = 1050.8571
$ws.Cells.Item(107, 14).Value = -5414.25
# row 122 (Leve Item ID 36182)
$ws.Cells.Item(122, 8).Value = 203599.8
$ws.Cells.Item(122, 10).Value = 253249.75
$ws.Cells.Item(122, 12).Value = 759749.25
$ws.Cells.Item(122, 14).Value = -764649.25

$ws = $wb.Worksheets.Item("LTW")
# row 5 (Leve Item ID 3790)
$ws.Cells.Item(5, 8).Value = 15208
$ws.Cells.Item(5, 9).Value = 8000
$ws.Cells.Item(5, 11).Value = 8000
$ws.Cells.Item(5, 13).Value = -7887
# row 16 (Leve Item ID 5289)
$ws.Cells.Item(16, 8).Value = 10499.5
$ws.Cells.Item(16, 9).Value = 10499.5
$ws.Cells.Item(16, 11).Value = 10499.5
$ws.Cells.Item(16, 13).Value = -10329.5
# row 22 (Leve Item ID 5277)
$ws.Cells.Item(22, 8).Value = 1000.5
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 1000.5
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 1000.5
$ws.Cells.Item(22, 13).Value = $null
$ws.Cells.Item(22, 14).Value = -1590.5
# row 27 (Leve Item ID 5277)
$ws.Cells.Item(27, 8).Value = 1000.5
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 1000.5
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 1000.5
$ws.Cells.Item(27, 13).Value = $null
$ws.Cells.Item(27, 14).Value = -1214.5
# row 40 (Leve Item ID 36248)
$ws.Cells.Item(40, 8).Value = 7834.2856
$ws.Cells.Item(40, 9).Value = 8488.5
$ws.Cells.Item(40, 10).Value = 6198.75
$ws.Cells.Item(40, 11).Value = 8488.5
$ws.Cells.Item(40, 12).Value = 6198.75
$ws.Cells.Item(40, 13).Value = -8352.5
$ws.Cells.Item(40, 14).Value = -6470.75
# row 46 (Leve Item ID 5282)
$ws.Cells.Item(46, 8).Value = 2360
$ws.Cells.Item(46, 9).Value = 1995
$ws.Cells.Item(46, 10).Value = 2725
$ws.Cells.Item(46, 11).Value = 1995
$ws.Cells.Item(46, 12).Value = 2725
$ws.Cells.Item(46, 13).Value = -1807
$ws.Cells.Item(46, 14).Value = -3101
# row 68 (Leve Item ID 12563)
$ws.Cells.Item(68, 8).Value = 37999.6
$ws.Cells.Item(68, 9).Value = 34999.5
$ws.Cells.Item(68, 11).Value = 34999.5
$ws.Cells.Item(68, 13).Value = -34250.5
# row 71 (Leve Item ID 12563)
$ws.Cells.Item(71, 8).Value = 37999.6
$ws.Cells.Item(71, 9).Value = 34999.5
$ws.Cells.Item(71, 11).Value = 174997.5
$ws.Cells.Item(71, 13).Value = -171253.5
# row 101 (Leve Item ID 18549)
$ws.Cells.Item(101, 8).Value = 12486
$ws.Cells.Item(101, 10).Value = 12486
$ws.Cells.Item(101, 12).Value = 12486
$ws.Cells.Item(101, 14).Value = -18976
# row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 4641.1113
$ws.Cells.Item(132, 9).Value = 4641.1113
$ws.Cells.Item(132, 11).Value = 13923.3339
$ws.Cells.Item(132, 13).Value = -11393.3339

$ws = $wb.Worksheets.Item("WVR")
# row 69 (Leve Item ID 10951)
$ws.Cells.Item(69, 8).Value = 8285.286
$ws.Cells.Item(69, 10).Value = 8285.286
$ws.Cells.Item(69, 12).Value = 8285.286
$ws.Cells.Item(69, 14).Value = -9783.286
# row 72 (Leve Item ID 10951)
$ws.Cells.Item(72, 8).Value = 8285.286
$ws.Cells.Item(72, 10).Value = 8285.286
$ws.Cells.Item(72, 12).Value = 24855.858
$ws.Cells.Item(72, 14).Value = -32343.858
# row 100 (Leve Item ID 19981)
$ws.Cells.Item(100, 8).Value = 2124.75
$ws.Cells.Item(100, 9).Value = 3250
$ws.Cells.Item(100, 10).Value = 999.5
$ws.Cells.Item(100, 11).Value = 6500
$ws.Cells.Item(100, 12).Value = 1999
$ws.Cells.Item(100, 13).Value = -5959
$ws.Cells.Item(100, 14).Value = -3081
# row 103 (Leve Item ID 18548)
$ws.Cells.Item(103, 8).Value = 19019
$ws.Cells.Item(103, 10).Value = 19019
$ws.Cells.Item(103, 12).Value = 19019
$ws.Cells.Item(103, 14).Value = -21363
# row 107 (Leve Item ID 27746)
$ws.Cells.Item(107, 8).Value = 2374.625
$ws.Cells.Item(107, 9).Value = 2416.3333
$ws.Cells.Item(107, 11).Value = 7248.999899999999
$ws.Cells.Item(107, 13).Value = -5328.999899999999
# row 112 (Leve Item ID 25836)
$ws.Cells.Item(112, 8).Value = 36128.832
$ws.Cells.Item(112, 10).Value = 36128.832
$ws.Cells.Item(112, 12).Value = 36128.832
$ws.Cells.Item(112, 14).Value = -39082.832
# row 122 (Leve Item ID 36208)
$ws.Cells.Item(122, 8).Value = 1444.3334
$ws.Cells.Item(122, 9).Value = 1447.25
$ws.Cells.Item(122, 11).Value = 4341.75
$ws.Cells.Item(122, 13).Value = -1891.75
# row 125 (Leve Item ID 34276)
$ws.Cells.Item(125, 8).Value = 98071.664
$ws.Cells.Item(125, 10).Value = 98071.664
$ws.Cells.Item(125, 12).Value = 98071.664
$ws.Cells.Item(125, 14).Value = -107911.664
